# Auto-generated: update leve-profit calculation columns (H,I,J,K,L,M,N)
# across all 8 sheets, per the scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 771.1111
$ws.Range("I4").Value = 156.83333
$ws.Range("K4").Value = 156.83333
$ws.Range("M4").Value = -42.83332999999999
$ws.Range("H40").Value = 2814.2856
$ws.Range("J40").Value = 1999
$ws.Range("L40").Value = 1999
$ws.Range("N40").Value = -2349
$ws.Range("H87").Value = 39998.332
$ws.Range("J87").Value = 39998.332
$ws.Range("L87").Value = 39998.332
$ws.Range("N87").Value = -42494.332
$ws.Range("H88").Value = 799.8387
$ws.Range("I88").Value = 1084.2
$ws.Range("J88").Value = 664.4286
$ws.Range("K88").Value = 1084.2
$ws.Range("L88").Value = 664.4286
$ws.Range("M88").Value = -678.2
$ws.Range("N88").Value = -1476.4286
$ws.Range("H90").Value = 39998.332
$ws.Range("J90").Value = 39998.332
$ws.Range("L90").Value = 119994.996
$ws.Range("N90").Value = -132474.996
$ws.Range("H91").Value = 799.8387
$ws.Range("I91").Value = 1084.2
$ws.Range("J91").Value = 664.4286
$ws.Range("K91").Value = 1084.2
$ws.Range("L91").Value = 664.4286
$ws.Range("M91").Value = 319.8
$ws.Range("N91").Value = -3472.4286
$ws.Range("H92").Value = 24089.143
$ws.Range("I92").Value = 29687.176
$ws.Range("J92").Value = 297.5
$ws.Range("K92").Value = 29687.176
$ws.Range("L92").Value = 297.5
$ws.Range("M92").Value = -28439.176
$ws.Range("N92").Value = -2793.5
$ws.Range("H94").Value = 1039.091
$ws.Range("I94").Value = 1039.091
$ws.Range("K94").Value = 1039.091
$ws.Range("M94").Value = -588.0909999999999
$ws.Range("H98").Value = 1541
$ws.Range("I98").Value = 1620.7273
$ws.Range("J98").Value = 1248.6666
$ws.Range("K98").Value = 1620.7273
$ws.Range("L98").Value = 1248.6666
$ws.Range("M98").Value = -122.7273
$ws.Range("N98").Value = -4244.6666
$ws.Range("H101").Value = 1793.625
$ws.Range("I101").Value = 1502.6666
$ws.Range("K101").Value = 4507.9998
$ws.Range("M101").Value = -2885.9998
$ws.Range("H122").Value = 1541
$ws.Range("I122").Value = 1620.7273
$ws.Range("J122").Value = 1248.6666
$ws.Range("K122").Value = 4862.1819
$ws.Range("L122").Value = 3745.9998
$ws.Range("M122").Value = -2412.1819
$ws.Range("N122").Value = -8645.9998
$ws.Range("H138").Value = 3366.4285
$ws.Range("J138").Value = 3511.068
$ws.Range("L138").Value = 10533.204
$ws.Range("N138").Value = -20813.204
$ws.Range("H141").Value = 3254.8333
$ws.Range("I141").Value = 3254.8333
$ws.Range("K141").Value = 9764.499899999999
$ws.Range("M141").Value = -4584.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31608
$ws.Range("H60").Value = 18202.334
$ws.Range("I60").Value = 9050
$ws.Range("J60").Value = 22778.5
$ws.Range("K60").Value = 9050
$ws.Range("L60").Value = 22778.5
$ws.Range("M60").Value = -8317
$ws.Range("N60").Value = -24244.5
$ws.Range("H76").Value = 188890.67
$ws.Range("J76").Value = 188890.67
$ws.Range("L76").Value = 188890.67
$ws.Range("N76").Value = -189566.67
$ws.Range("H79").Value = 188890.67
$ws.Range("J79").Value = 188890.67
$ws.Range("L79").Value = 188890.67
$ws.Range("N79").Value = -191230.67
$ws.Range("H103").Value = 78500
$ws.Range("J103").Value = 78500
$ws.Range("L103").Value = 78500
$ws.Range("N103").Value = -80844
$ws.Range("H132").Value = 8536.581
$ws.Range("I132").Value = 8771.536
$ws.Range("K132").Value = 26314.608
$ws.Range("M132").Value = -23784.608
$ws.Range("H133").Value = 79333.336
$ws.Range("J133").Value = 79333.336
$ws.Range("L133").Value = 79333.336
$ws.Range("N133").Value = -84393.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3004
$ws.Range("I5").Value = 3004
$ws.Range("K5").Value = 3004
$ws.Range("M5").Value = -2891
$ws.Range("H94").Value = 101719.664
$ws.Range("I94").Value = 996.1667
$ws.Range("K94").Value = 996.1667
$ws.Range("M94").Value = -545.1667
$ws.Range("H134").Value = 5618.1665
$ws.Range("I134").Value = 5132.0625
$ws.Range("K134").Value = 15396.1875
$ws.Range("M134").Value = -12861.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1764.8
$ws.Range("I10").Value = 1764.8
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1764.8
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = -1625.8
$ws.Range("H31").Value = 1882.9584
$ws.Range("I31").Value = 1234.091
$ws.Range("J31").Value = 3310.4666
$ws.Range("K31").Value = 1234.091
$ws.Range("L31").Value = 3310.4666
$ws.Range("M31").Value = -939.0909999999999
$ws.Range("N31").Value = -3900.4666
$ws.Range("H34").Value = 1882.9584
$ws.Range("I34").Value = 1234.091
$ws.Range("J34").Value = 3310.4666
$ws.Range("K34").Value = 1234.091
$ws.Range("L34").Value = 3310.4666
$ws.Range("M34").Value = -1032.091
$ws.Range("N34").Value = -3714.4666
$ws.Range("H41").Value = 46554.668
$ws.Range("J41").Value = 49999
$ws.Range("L41").Value = 49999
$ws.Range("N41").Value = -50855
$ws.Range("H74").Value = 73918.82000000001
$ws.Range("J74").Value = 73918.82000000001
$ws.Range("L74").Value = 73918.82000000001
$ws.Range("N74").Value = -75666.82000000001
$ws.Range("H77").Value = 73918.82000000001
$ws.Range("J77").Value = 73918.82000000001
$ws.Range("L77").Value = 221756.46
$ws.Range("N77").Value = -230492.46

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 895.0909
$ws.Range("I5").Value = 862.3333
$ws.Range("J5").Value = 907.375
$ws.Range("K5").Value = 2586.9999
$ws.Range("L5").Value = 2722.125
$ws.Range("M5").Value = -2474.9999
$ws.Range("N5").Value = -2946.125
$ws.Range("H109").Value = 1269.8
$ws.Range("I109").Value = 1269.8
$ws.Range("K109").Value = 3809.4
$ws.Range("M109").Value = -2769.4
$ws.Range("H114").Value = 663.3333
$ws.Range("I114").Value = 490
$ws.Range("K114").Value = 1470
$ws.Range("M114").Value = 1784
$ws.Range("H135").Value = 895.0909
$ws.Range("I135").Value = 862.3333
$ws.Range("J135").Value = 907.375
$ws.Range("K135").Value = 7760.9997
$ws.Range("L135").Value = 8166.375
$ws.Range("M135").Value = -5225.9997
$ws.Range("N135").Value = -13236.375
$ws.Range("H140").Value = 2507.2
$ws.Range("I140").Value = 2285.7778
$ws.Range("K140").Value = 6857.3334
$ws.Range("M140").Value = -1677.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 47500
$ws.Range("I5").Value = 47500
$ws.Range("K5").Value = 47500
$ws.Range("M5").Value = -47388
$ws.Range("H9").Value = 3538.6
$ws.Range("I9").Value = 3538.6
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 3538.6
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = -3368.6
$ws.Range("H43").Value = 178169.33
$ws.Range("I43").Value = 178169.33
$ws.Range("K43").Value = 178169.33
$ws.Range("M43").Value = -178018.33
$ws.Range("H45").Value = 48633
$ws.Range("J45").Value = 48633
$ws.Range("L45").Value = 48633
$ws.Range("N45").Value = -49751

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 81141.42999999999
$ws.Range("J2").Value = 11331.667
$ws.Range("L2").Value = 11331.667
$ws.Range("N2").Value = -11555.667
$ws.Range("H12").Value = 3371.75
$ws.Range("J12").Value = 3371.75
$ws.Range("L12").Value = 3371.75
$ws.Range("N12").Value = -3711.75
$ws.Range("H22").Value = 3468.9
$ws.Range("I22").Value = 2916.3333
$ws.Range("J22").Value = 3705.7144
$ws.Range("K22").Value = 2916.3333
$ws.Range("L22").Value = 3705.7144
$ws.Range("M22").Value = -2621.3333
$ws.Range("N22").Value = -4295.7144
$ws.Range("H25").Value = 7386
$ws.Range("I25").Value = 3249.5
$ws.Range("J25").Value = 9040.6
$ws.Range("K25").Value = 3249.5
$ws.Range("L25").Value = 9040.6
$ws.Range("M25").Value = -3019.5
$ws.Range("N25").Value = -9500.6
$ws.Range("H27").Value = 3468.9
$ws.Range("I27").Value = 2916.3333
$ws.Range("J27").Value = 3705.7144
$ws.Range("K27").Value = 2916.3333
$ws.Range("L27").Value = 3705.7144
$ws.Range("M27").Value = -2809.3333
$ws.Range("N27").Value = -3919.7144
$ws.Range("H93").Value = 2097.3076
$ws.Range("I93").Value = 1541.4445
$ws.Range("J93").Value = 3348
$ws.Range("K93").Value = 1541.4445
$ws.Range("L93").Value = 3348
$ws.Range("M93").Value = -293.4445000000001
$ws.Range("N93").Value = -5844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8499.666999999999
$ws.Range("I2").Value = 4999.3335
$ws.Range("J2").Value = 12000
$ws.Range("K2").Value = 4999.3335
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = -4887.3335
$ws.Range("N2").Value = -12224
